$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4315
$ws.Range("K3").Value = 4381
$ws.Range("J4").Value = 1008
$ws.Range("K4").Value = 887
$ws.Range("K5").Value = 322
$ws.Range("K6").Value = 4930
$ws.Range("J7").Value = 14625
$ws.Range("K7").Value = 14835

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 433
$ws.Range("K8").Value = 996
$ws.Range("K9").Value = 62
$ws.Range("K10").Value = 81
$ws.Range("K11").Value = 295
$ws.Range("K12").Value = 26
$ws.Range("K19").Value = 451
$ws.Range("K20").Value = 336
$ws.Range("K23").Value = 152
$ws.Range("K29").Value = 788
$ws.Range("K30").Value = 54
$ws.Range("K33").Value = 621
$ws.Range("K36").Value = 187
$ws.Range("K37").Value = 506
$ws.Range("K42").Value = 544
$ws.Range("K43").Value = 131
$ws.Range("K44").Value = 134
$ws.Range("K47").Value = 93
$ws.Range("K49").Value = 87
$ws.Range("K51").Value = 186
$ws.Range("K54").Value = 279
$ws.Range("K55").Value = 167
$ws.Range("J63").Value = 73
$ws.Range("K63").Value = 53
$ws.Range("K64").Value = 92
$ws.Range("K65").Value = 338
$ws.Range("K67").Value = 571
$ws.Range("K69").Value = 32
$ws.Range("K72").Value = 68
$ws.Range("K79").Value = 380
$ws.Range("K84").Value = 108
$ws.Range("K85").Value = 666
$ws.Range("K88").Value = 170
$ws.Range("K89").Value = 206
$ws.Range("K90").Value = 133
$ws.Range("K95").Value = 261
$ws.Range("K96").Value = 163
$ws.Range("K97").Value = 126
$ws.Range("K99").Value = 247
$ws.Range("J101").Value = 14625
$ws.Range("K101").Value = 14835

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 140
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 433

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 96
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 295

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 206

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 236
$ws.Range("K3").Value = 221
$ws.Range("K7").Value = 666

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value = 299
$ws.Range("K6").Value = 335
$ws.Range("K7").Value = 996

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 231
$ws.Range("K4").Value = 27
$ws.Range("K6").Value = 180
$ws.Range("K7").Value = 621

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 90
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 140
$ws.Range("K7").Value = 506

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 103
$ws.Range("K7").Value = 338

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 66
$ws.Range("K7").Value = 247

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 199
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 571

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 77
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 226
$ws.Range("K3").Value = 279
$ws.Range("K6").Value = 219
$ws.Range("K7").Value = 788

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 143
$ws.Range("K3").Value = 136
$ws.Range("K6").Value = 139
$ws.Range("K7").Value = 451

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 151
$ws.Range("K7").Value = 544

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 39
$ws.Range("K3").Value = 77

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 380

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 336

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 74
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 53
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 186

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 27
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 26
